# Updates the cryptocurrency price table (columns B-E) on the active sheet
# to reflect the latest scraped values, as produced by the GitHub Actions job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; B = $null; C = $null; D = "30.607.96"; E = "  +2.10%  " },
    @{ Row = 3; B = $null; C = $null; D = "1.887.60"; E = $null },
    @{ Row = 4; B = $null; C = $null; D = "1.001"; E = "  +0.04%  " },
    @{ Row = 5; B = $null; C = $null; D = "245.11"; E = "  +1.02%  " },
    @{ Row = 6; B = $null; C = $null; D = "1.000"; E = "  +0.04%  " },
    @{ Row = 7; B = $null; C = $null; D = "0.4911"; E = "  -0.13%  " },
    @{ Row = 8; B = $null; C = $null; D = "0.2948"; E = "  +0.28%  " },
    @{ Row = 9; B = $null; C = $null; D = "0.06773"; E = "  +2.39%  " },
    @{ Row = 10; B = $null; C = $null; D = "1.888.05"; E = "  +0.35%  " },
    @{ Row = 11; B = $null; C = $null; D = "17.27"; E = "  +3.63%  " },
    @{ Row = 12; B = $null; C = $null; D = "0.07237"; E = "  +0.93%  " },
    @{ Row = 13; B = $null; C = $null; D = "91.16"; E = "  +5.61%  " },
    @{ Row = 14; B = $null; C = $null; D = "0.6775"; E = "  +1.72%  " },
    @{ Row = 15; B = $null; C = $null; D = "5.045"; E = "  +3.55%  " },
    @{ Row = 16; B = $null; C = $null; D = "30.600.99"; E = "  +2.08%  " },
    @{ Row = 17; B = $null; C = $null; D = "0.000007961"; E = "  +1.96%  " },
    @{ Row = 18; B = $null; C = $null; D = "1.0000"; E = "  +0.02%  " },
    @{ Row = 19; B = $null; C = $null; D = "13.14"; E = "  +2.80%  " },
    @{ Row = 20; B = $null; C = $null; D = "2.131.65"; E = "  +0.51%  " },
    @{ Row = 21; B = $null; C = $null; D = $null; E = "  +0.02%  " },
    @{ Row = 22; B = $null; C = $null; D = "4.823"; E = "  +0.99%  " },
    @{ Row = 23; B = $null; C = $null; D = "192.32"; E = "  +36.75%  " },
    @{ Row = 24; B = $null; C = $null; D = "6.068"; E = "  +3.84%  " },
    @{ Row = 25; B = $null; C = $null; D = "9.323"; E = "  +2.68%  " },
    @{ Row = 26; B = $null; C = $null; D = "155.56"; E = "  +3.35%  " },
    @{ Row = 27; B = $null; C = $null; D = "19.14"; E = "  +12.99%  " },
    @{ Row = 28; B = $null; C = $null; D = "1.903"; E = "  +0.19%  " },
    @{ Row = 29; B = $null; C = $null; D = "1.400"; E = "  +0.74%  " },
    @{ Row = 30; B = $null; C = $null; D = "4.323"; E = "  +3.28%  " },
    @{ Row = 31; B = $null; C = $null; D = "0.09034"; E = "  +3.28%  " },
    @{ Row = 32; B = $null; C = $null; D = "4.007"; E = "  +0.70%  " },
    @{ Row = 33; B = $null; C = $null; D = "0.05199"; E = "  +3.67%  " },
    @{ Row = 34; B = $null; C = $null; D = "0.7536"; E = "  +5.25%  " },
    @{ Row = 35; B = $null; C = $null; D = "1.110"; E = "  +0.09%  " },
    @{ Row = 36; B = $null; C = $null; D = "2.761"; E = "  +3.44%  " },
    @{ Row = 37; B = $null; C = $null; D = "0.01833"; E = "  +2.49%  " },
    @{ Row = 38; B = $null; C = $null; D = "2.670"; E = "  -0.82%  " },
    @{ Row = 39; B = $null; C = $null; D = "2.142"; E = "  -0.72%  " },
    @{ Row = 40; B = $null; C = $null; D = "0.9327"; E = "  -0.74%  " },
    @{ Row = 41; B = $null; C = $null; D = "0.4414"; E = "  +4.47%  " },
    @{ Row = 42; B = $null; C = $null; D = "105.11"; E = "  +1.37%  " },
    @{ Row = 43; B = $null; C = $null; D = $null; E = "  +0.09%  " },
    @{ Row = 44; B = $null; C = $null; D = "5.732"; E = "  -0.14%  " },
    @{ Row = 45; B = $null; C = $null; D = "7.588"; E = "  +3.68%  " },
    @{ Row = 46; B = $null; C = $null; D = $null; E = "  +5.84%  " },
    @{ Row = 47; B = $null; C = $null; D = $null; E = "  +2.73%  " },
    @{ Row = 48; B = "NEARProtocol"; C = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D = "1.434"; E = "  +7.20%  " },
    @{ Row = 49; B = "EnergySwap"; C = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D = "8.708"; E = "  +6.01%  " },
    @{ Row = 50; B = $null; C = $null; D = "0.3918"; E = "  +4.42%  " },
    @{ Row = 51; B = $null; C = $null; D = "33.54"; E = "  +2.67%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.B) {
        $ws.Range("B$row").Value = $u.B
    }
    if ($null -ne $u.C) {
        $ws.Range("C$row").Value = $u.C
    }
    if ($null -ne $u.D) {
        # The Price column holds values such as "1.000" or "30.607.96" that
        # Excel would otherwise auto-convert to numbers, losing the original
        # formatting. Force the cell to text so the literal string is kept,
        # then restore the default (Normal) style so no extra formatting
        # lingers on the cell.
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $ws.Range("E$row").Value = $u.E
    }
}
